$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("ECs", "Icosl", "Cd28", "M2", 3, 1, 0.9626043333333333, 2.887813, 0.07904870764102041, 0.07904870764102041, 3, 1, 7.543204, 22.629612, 0.9718881576768906, 0.9718881576768906, 7.261120857617334, 65.350087718556, 0.07682650283597048, 0.07682650283597048),
    @("ECs", "Icosl", "Cd28", "sCs", 3, 1, 0.9626043333333333, 2.887813, 0.07904870764102041, 0.07904870764102041, 1, 0.3333333333333333, 0.218187, 0.6545609999999999, 0.02811184232310935, 0.02811184232310934, 0.210027751677, 1.890249765093, 0.002222204805049935, 0.002222204805049934),
    @("FAPs", "Icosl", "Cd28", "M2", 3, 1, 2.399281333333333, 7.197844, 0.1970280852678733, 0.1970280852678733, 3, 1, 7.543204, 22.629612, 0.9718881576768906, 0.9718881576768906, 18.09826855072533, 162.884416956528, 0.1914892628015987, 0.1914892628015987),
    @("FAPs", "Icosl", "Cd28", "sCs", 3, 1, 2.399281333333333, 7.197844, 0.1970280852678733, 0.1970280852678733, 1, 0.3333333333333333, 0.218187, 0.6545609999999999, 0.02811184232310935, 0.02811184232310934, 0.523491996276, 4.711427966484, 0.005538822466274597, 0.005538822466274597),
    @("M2", "Icosl", "Cd28", "M2", 3, 1, 7.460122666666666, 22.380368, 0.6126224817640369, 0.6126224817640369, 3, 1, 7.543204, 22.629612, 0.9718881576768906, 0.9718881576768906, 56.27322713969066, 506.459044257216, 0.5954005351530943, 0.5954005351530943),
    @("M2", "Icosl", "Cd28", "sCs", 3, 1, 7.460122666666666, 22.380368, 0.6126224817640369, 0.6126224817640369, 1, 0.3333333333333333, 0.218187, 0.6545609999999999, 0.02811184232310935, 0.02811184232310934, 1.627701784272, 14.649316058448, 0.01722194661094253, 0.01722194661094253),
    @("sCs", "Icosl", "Cd28", "M2", 3, 1, 1.355348666666667, 4.066046, 0.1113007253270695, 0.1113007253270695, 3, 1, 7.543204, 22.629612, 0.9718881576768906, 0.9718881576768906, 10.22367148379467, 92.013043354152, 0.1081718568862272, 0.1081718568862272),
    @("sCs", "Icosl", "Cd28", "sCs", 3, 1, 1.355348666666667, 4.066046, 0.1113007253270695, 0.1113007253270695, 1, 0.3333333333333333, 0.218187, 0.6545609999999999, 0.02811184232310935, 0.02811184232310934, 0.295719459534, 2.661475135806, 0.00312886844084228, 0.00312886844084228),
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $rowNum = $i + 2
    $rowData = $data[$i]
    for ($j = 0; $j -lt $rowData.Count; $j++) {
        $colNum = $j + 1
        $ws.Cells.Item($rowNum, $colNum).Value = $rowData[$j]
    }
}
